$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prefix with an apostrophe so Excel stores the date-like text literally
# (as a string) instead of auto-converting it to a date serial number,
# then reset the style to Normal so no stray number-format/quote-prefix
# style is left attached to the cell.
$ws.Range("A45").Value = "'2025-09-29"
$ws.Range("A45").Style = "Normal"

$ws.Range("B45").Value = 55.27999877929688
$ws.Range("C45").Value = 672.5
$ws.Range("D45").Value = 324.8500061035156
